$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H9").Value = 144
$ws.Range("J9").Value = 210
$ws.Range("L9").Value = 210
$ws.Range("N9").Value = -548
$ws.Range("H12").Value = 546.9231
$ws.Range("I12").Value = 514.1111
$ws.Range("J12").Value = 620.75
$ws.Range("K12").Value = 514.1111
$ws.Range("L12").Value = 620.75
$ws.Range("M12").Value = -344.1111
$ws.Range("N12").Value = -960.75
$ws.Range("H13").Value = 2105
$ws.Range("I13").Value = 2105
$ws.Range("J13").Value = 0
$ws.Range("K13").Value = 2105
$ws.Range("L13").Value = 0
$ws.Range("M13").Value = -1936
$ws.Range("N13").ClearContents()
$ws.Range("H15").Value = 1708.1666
$ws.Range("I15").Value = 1708.1666
$ws.Range("K15").Value = 5124.4998
$ws.Range("M15").Value = -4955.4998
$ws.Range("H75").Value = 22267.5
$ws.Range("I75").Value = 14285
$ws.Range("J75").Value = 30250
$ws.Range("K75").Value = 14285
$ws.Range("L75").Value = 30250
$ws.Range("M75").Value = -13349
$ws.Range("N75").Value = -32122
$ws.Range("H78").Value = 22267.5
$ws.Range("I78").Value = 14285
$ws.Range("J78").Value = 30250
$ws.Range("K78").Value = 42855
$ws.Range("L78").Value = 90750
$ws.Range("M78").Value = -38175
$ws.Range("N78").Value = -100110
$ws.Range("H86").Value = 1840
$ws.Range("H89").Value = 1840
$ws.Range("H92").Value = 1528.6666
$ws.Range("J92").Value = 1999.5
$ws.Range("L92").Value = 1999.5
$ws.Range("N92").Value = -4495.5
$ws.Range("H98").Value = 441.2
$ws.Range("I98").Value = 387.07144
$ws.Range("J98").Value = 1199
$ws.Range("K98").Value = 387.07144
$ws.Range("L98").Value = 1199
$ws.Range("M98").Value = 1110.92856
$ws.Range("N98").Value = -4195
$ws.Range("H106").Value = 10000
$ws.Range("I106").Value = 10000
$ws.Range("K106").Value = 10000
$ws.Range("M106").Value = -9369
$ws.Range("H107").Value = 1596.75
$ws.Range("I107").Value = 1883.8889
$ws.Range("J107").Value = 735.3333
$ws.Range("K107").Value = 1883.8889
$ws.Range("L107").Value = 735.3333
$ws.Range("M107").Value = 36.11110000000008
$ws.Range("N107").Value = -4575.3333
$ws.Range("H116").Value = 4389.5
$ws.Range("I116").Value = 4389.5
$ws.Range("K116").Value = 4389.5
$ws.Range("M116").Value = -947.5
$ws.Range("H122").Value = 441.2
$ws.Range("I122").Value = 387.07144
$ws.Range("J122").Value = 1199
$ws.Range("K122").Value = 1161.21432
$ws.Range("L122").Value = 3597
$ws.Range("M122").Value = 1288.78568
$ws.Range("N122").Value = -8497

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H63").Value = 3720.2222
$ws.Range("I63").Value = 1685.625
$ws.Range("K63").Value = 1685.625
$ws.Range("M63").Value = -999.625
$ws.Range("H66").Value = 3720.2222
$ws.Range("I66").Value = 1685.625
$ws.Range("K66").Value = 8428.125
$ws.Range("M66").Value = -4996.125
$ws.Range("H74").Value = 2768.5
$ws.Range("I74").Value = 2768.5
$ws.Range("K74").Value = 2768.5
$ws.Range("M74").Value = -1894.5
$ws.Range("H77").Value = 2768.5
$ws.Range("I77").Value = 2768.5
$ws.Range("K77").Value = 13842.5
$ws.Range("M77").Value = -9474.5
$ws.Range("H88").Value = 1378.5
$ws.Range("I88").Value = 1443.8
$ws.Range("J88").Value = 1269.6666
$ws.Range("K88").Value = 1443.8
$ws.Range("L88").Value = 1269.6666
$ws.Range("M88").Value = -1037.8
$ws.Range("N88").Value = -2081.6666
$ws.Range("H91").Value = 1378.5
$ws.Range("I91").Value = 1443.8
$ws.Range("J91").Value = 1269.6666
$ws.Range("K91").Value = 1443.8
$ws.Range("L91").Value = 1269.6666
$ws.Range("M91").Value = -39.79999999999995
$ws.Range("N91").Value = -4077.6666

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H22").Value = 449.85715
$ws.Range("I22").Value = 274.83334
$ws.Range("K22").Value = 274.83334
$ws.Range("M22").Value = -101.83334
$ws.Range("H107").Value = 10000
$ws.Range("I107").Value = 0
$ws.Range("K107").Value = 0
$ws.Range("M107").ClearContents()

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 2702.8667
$ws.Range("I58").Value = 1164.9
$ws.Range("K58").Value = 1164.9
$ws.Range("M58").Value = -961.9000000000001
$ws.Range("H105").Value = 1991
$ws.Range("I105").Value = 1947.8334
$ws.Range("J105").Value = 2250
$ws.Range("K105").Value = 1947.8334
$ws.Range("L105").Value = 2250
$ws.Range("M105").Value = -200.8334
$ws.Range("N105").Value = -5744
$ws.Range("H132").Value = 1935.6923
$ws.Range("I132").Value = 1935.6923
$ws.Range("K132").Value = 5807.0769
$ws.Range("M132").Value = -3277.0769
$ws.Range("H136").Value = 2702.8667
$ws.Range("I136").Value = 1164.9
$ws.Range("K136").Value = 3494.7
$ws.Range("M136").Value = -944.7000000000003

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H126").Value = 2451.4285
$ws.Range("I126").Value = 1908.8334
$ws.Range("J126").Value = 5707
$ws.Range("K126").Value = 5726.5002
$ws.Range("L126").Value = 17121
$ws.Range("M126").Value = -3256.5002
$ws.Range("N126").Value = -22061

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 7929.75
$ws.Range("I7").Value = 7719.857
$ws.Range("J7").Value = 9399
$ws.Range("K7").Value = 7719.857
$ws.Range("L7").Value = 9399
$ws.Range("M7").Value = -7607.857
$ws.Range("N7").Value = -9623
$ws.Range("H48").Value = 2041
$ws.Range("I48").Value = 2041
$ws.Range("K48").Value = 2041
$ws.Range("M48").Value = -1380
$ws.Range("H126").Value = 7929.75
$ws.Range("I126").Value = 7719.857
$ws.Range("J126").Value = 9399
$ws.Range("K126").Value = 23159.571
$ws.Range("L126").Value = 28197
$ws.Range("M126").Value = -20689.571
$ws.Range("N126").Value = -33137

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 5829
$ws.Range("H65").Value = 5829
$ws.Range("H132").Value = 1159.0358
$ws.Range("I132").Value = 1159.0358
$ws.Range("K132").Value = 3477.1074
$ws.Range("M132").Value = -947.1074000000003
